# Weekly market-price refresh for the Excalibur Profits leve-crafting workbook.
# Updates currentAveragePrice / HQ / NQ columns (H-N) per leve row from the
# latest Universalis snapshot; values below mirror the scheduled runner output.
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets("ALC")
# row 17
$ws.Range("H17").Value = 1759.1
$ws.Range("J17").Value = 1773.4359
$ws.Range("L17").Value = 5320.307699999999
$ws.Range("N17").Value = -5656.307699999999
# row 53
$ws.Range("H53").Value = 204.77777
$ws.Range("I53").Value = 97.25
$ws.Range("K53").Value = 97.25
$ws.Range("M53").Value = 539.75
# row 62
$ws.Range("H62").Value = 19412.584
$ws.Range("I62").Value = 17910.166
$ws.Range("J62").Value = 20915
$ws.Range("K62").Value = 17910.166
$ws.Range("L62").Value = 20915
$ws.Range("M62").Value = -17286.166
$ws.Range("N62").Value = -22163
# row 65
$ws.Range("H65").Value = 19412.584
$ws.Range("I65").Value = 17910.166
$ws.Range("J65").Value = 20915
$ws.Range("K65").Value = 89550.83
$ws.Range("L65").Value = 104575
$ws.Range("M65").Value = -86430.83
$ws.Range("N65").Value = -110815
# row 74
$ws.Range("H74").Value = 7601.4287
$ws.Range("J74").Value = 8773.267
$ws.Range("L74").Value = 8773.267
$ws.Range("N74").Value = -10645.267
# row 77
$ws.Range("H77").Value = 7601.4287
$ws.Range("J77").Value = 8773.267
$ws.Range("L77").Value = 43866.335
$ws.Range("N77").Value = -53226.335
# row 106
$ws.Range("H106").Value = 4068.3333
$ws.Range("I106").Value = 3602.5
$ws.Range("J106").Value = 5000
$ws.Range("K106").Value = 3602.5
$ws.Range("L106").Value = 5000
$ws.Range("M106").Value = -2971.5
$ws.Range("N106").Value = -6262
# row 123
$ws.Range("H123").Value = 4229954.5
$ws.Range("J123").Value = 75945.60000000001
$ws.Range("L123").Value = 75945.60000000001
$ws.Range("N123").Value = -85745.60000000001
# row 125
$ws.Range("H125").Value = 3024.647
$ws.Range("I125").Value = 3440.5715
$ws.Range("J125").Value = 2733.5
$ws.Range("K125").Value = 30965.1435
$ws.Range("L125").Value = 24601.5
$ws.Range("M125").Value = -28505.1435
$ws.Range("N125").Value = -29521.5
# row 132
$ws.Range("H132").Value = 43944.438
$ws.Range("I132").Value = 45557.66
$ws.Range("K132").Value = 136672.98
$ws.Range("M132").Value = -134142.98
# row 135
$ws.Range("H135").Value = 1398.0312
$ws.Range("I135").Value = 1410.871
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 12697.839
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -10162.839
$ws.Range("N135").Value = -14070
# row 137
$ws.Range("H137").Value = 1277208.5
$ws.Range("I137").Value = 1069697.4
$ws.Range("K137").Value = 3209092.2
$ws.Range("M137").Value = -3206542.2

# ----- ARM -----
$ws = $wb.Worksheets("ARM")
# row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# row 132
$ws.Range("H132").Value = 383473.6
$ws.Range("I132").Value = 448561.53
$ws.Range("K132").Value = 1345684.59
$ws.Range("M132").Value = -1343154.59

# ----- BSM -----
$ws = $wb.Worksheets("BSM")
# row 26
$ws.Range("H26").Value = 10234.5
$ws.Range("I26").Value = 10234.5
$ws.Range("K26").Value = 10234.5
$ws.Range("M26").Value = -9942.5
# row 94
$ws.Range("H94").Value = 1869.8422
$ws.Range("I94").Value = 1348.8572
$ws.Range("J94").Value = 3328.6
$ws.Range("K94").Value = 1348.8572
$ws.Range("L94").Value = 3328.6
$ws.Range("M94").Value = -897.8571999999999
$ws.Range("N94").Value = -4230.6
# row 105
$ws.Range("H105").Value = 1505.0667
$ws.Range("I105").Value = 1505.4286
$ws.Range("K105").Value = 1505.4286
$ws.Range("M105").Value = 241.5714
# row 107
$ws.Range("H107").Value = 6309.5
$ws.Range("I107").Value = 6119.6
$ws.Range("J107").Value = 6499.4
$ws.Range("K107").Value = 6119.6
$ws.Range("L107").Value = 6499.4
$ws.Range("M107").Value = -4199.6
$ws.Range("N107").Value = -10339.4
# row 134
$ws.Range("H134").Value = 390425.12
$ws.Range("I134").Value = 469358.16
$ws.Range("K134").Value = 1408074.48
$ws.Range("M134").Value = -1405539.48

# ----- CRP -----
$ws = $wb.Worksheets("CRP")
# row 16
$ws.Range("H16").Value = 639.8
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 500
$ws.Range("N16").Value = -1074
# row 58
$ws.Range("H58").Value = 1547775.8
$ws.Range("I58").Value = 3088298
$ws.Range("J58").Value = 7253.5
$ws.Range("K58").Value = 3088298
$ws.Range("L58").Value = 7253.5
$ws.Range("M58").Value = -3088095
$ws.Range("N58").Value = -7659.5
# row 113
$ws.Range("H113").Value = 639.8
$ws.Range("J113").Value = 500
$ws.Range("L113").Value = 500
$ws.Range("N113").Value = -4840
# row 136
$ws.Range("H136").Value = 1547775.8
$ws.Range("I136").Value = 3088298
$ws.Range("J136").Value = 7253.5
$ws.Range("K136").Value = 9264894
$ws.Range("L136").Value = 21760.5
$ws.Range("M136").Value = -9262344
$ws.Range("N136").Value = -26860.5
# row 137
$ws.Range("H137").Value = 63031.5

# ----- CUL -----
$ws = $wb.Worksheets("CUL")
# row 131
$ws.Range("H131").Value = 10049.5
$ws.Range("I131").Value = 674.25
$ws.Range("J131").Value = 14216.277
$ws.Range("K131").Value = 2022.75
$ws.Range("L131").Value = 42648.831
$ws.Range("M131").Value = 3017.25
$ws.Range("N131").Value = -52728.831
# row 132
$ws.Range("H132").Value = 3725.7273
$ws.Range("I132").Value = 2297.8
$ws.Range("J132").Value = 4915.6665
$ws.Range("K132").Value = 20680.2
$ws.Range("L132").Value = 44240.9985
$ws.Range("M132").Value = -18150.2
$ws.Range("N132").Value = -49300.9985
# row 136
$ws.Range("H136").Value = 6887
$ws.Range("I136").Value = 6887
$ws.Range("K136").Value = 20661
$ws.Range("M136").Value = -15561

# ----- GSM -----
$ws = $wb.Worksheets("GSM")
# row 19
$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 5000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -4712
$ws.Range("N19").ClearContents()
# row 97
$ws.Range("H97").Value = 2947.9285
$ws.Range("I97").Value = 796.2273
$ws.Range("K97").Value = 796.2273
$ws.Range("M97").Value = -300.2273
# row 113
$ws.Range("H113").Value = 3424.2856
$ws.Range("I113").Value = 2548.889
$ws.Range("K113").Value = 2548.889
$ws.Range("M113").Value = -378.8890000000001

# ----- LTW -----
$ws = $wb.Worksheets("LTW")
# row 35
$ws.Range("H35").Value = 3447.5
$ws.Range("I35").Value = 3447.5
$ws.Range("K35").Value = 3447.5
$ws.Range("M35").Value = -3111.5
# row 43
$ws.Range("H43").Value = 509333.34
$ws.Range("J43").Value = 610000
$ws.Range("L43").Value = 610000
$ws.Range("N43").Value = -610386
# row 61
$ws.Range("H61").Value = 4004.625
$ws.Range("I61").Value = 2113.3333
$ws.Range("J61").Value = 5139.4
$ws.Range("K61").Value = 2113.3333
$ws.Range("L61").Value = 5139.4
$ws.Range("M61").Value = -1911.3333
$ws.Range("N61").Value = -5543.4
# row 93
$ws.Range("H93").Value = 1853.8334
$ws.Range("I93").Value = 1639
$ws.Range("J93").Value = 2025.7
$ws.Range("K93").Value = 1639
$ws.Range("L93").Value = 2025.7
$ws.Range("M93").Value = -391
$ws.Range("N93").Value = -4521.7
# row 113
$ws.Range("H113").Value = 4004.625
$ws.Range("I113").Value = 2113.3333
$ws.Range("J113").Value = 5139.4
$ws.Range("K113").Value = 2113.3333
$ws.Range("L113").Value = 5139.4
$ws.Range("M113").Value = 56.66670000000022
$ws.Range("N113").Value = -9479.4
# row 122
$ws.Range("H122").Value = 3268.7292
$ws.Range("I122").Value = 3039.1936
$ws.Range("J122").Value = 3687.2942
$ws.Range("K122").Value = 9117.5808
$ws.Range("L122").Value = 11061.8826
$ws.Range("M122").Value = -6667.5808
$ws.Range("N122").Value = -15961.8826
# row 132
$ws.Range("H132").Value = 1237937.8
$ws.Range("I132").Value = 1443469
$ws.Range("K132").Value = 4330407
$ws.Range("M132").Value = -4327877
# row 135
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140
# row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ----- WVR -----
$ws = $wb.Worksheets("WVR")
# row 81
$ws.Range("H81").Value = 1583
$ws.Range("I81").Value = 1441.1428
$ws.Range("K81").Value = 2882.2856
$ws.Range("M81").Value = -1821.2856
# row 84
$ws.Range("H84").Value = 1583
$ws.Range("I84").Value = 1441.1428
$ws.Range("K84").Value = 14411.428
$ws.Range("M84").Value = -9107.428
# row 113
$ws.Range("H113").Value = 3012.28
$ws.Range("I113").Value = 1408.5
$ws.Range("J113").Value = 5053.4546
$ws.Range("K113").Value = 4225.5
$ws.Range("L113").Value = 15160.3638
$ws.Range("M113").Value = -2055.5
$ws.Range("N113").Value = -19500.3638
# row 118
$ws.Range("H118").Value = 110000
$ws.Range("J118").Value = 110000
$ws.Range("L118").Value = 110000
$ws.Range("N118").Value = -113314
# row 132
$ws.Range("H132").Value = 5752615
$ws.Range("I132").Value = 6942136
$ws.Range("K132").Value = 20826408
$ws.Range("M132").Value = -20823878
